$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (28 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 40001720
$ws.Range("I125").Value = 62501304
$ws.Range("J125").Value = 2461.2222
$ws.Range("K125").Value = 562511736
$ws.Range("L125").Value = 22150.9998
$ws.Range("M125").Value = -562509276
$ws.Range("N125").Value = -27070.9998
$ws.Range("H137").Value = 22388.117
$ws.Range("I137").Value = 31931.906
$ws.Range("J137").Value = 6314.3687
$ws.Range("K137").Value = 95795.71799999999
$ws.Range("L137").Value = 18943.1061
$ws.Range("M137").Value = -93245.71799999999
$ws.Range("N137").Value = -24043.1061
$ws.Range("H138").Value = 1459.35
$ws.Range("I138").Value = 740.76
$ws.Range("J138").Value = 2177.94
$ws.Range("K138").Value = 2222.28
$ws.Range("L138").Value = 6533.82
$ws.Range("M138").Value = 2917.72
$ws.Range("N138").Value = -16813.82
$ws.Range("H141").Value = 1468.4265
$ws.Range("I141").Value = 793.68085
$ws.Range("J141").Value = 2978.5715
$ws.Range("K141").Value = 2381.04255
$ws.Range("L141").Value = 8935.7145
$ws.Range("M141").Value = 2798.95745
$ws.Range("N141").Value = -19295.7145

# --- Sheet: ARM (32 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 872.63635
$ws.Range("I61").Value = 846.8946999999999
$ws.Range("J61").Value = 1035.6666
$ws.Range("K61").Value = 846.8946999999999
$ws.Range("L61").Value = 1035.6666
$ws.Range("M61").Value = -634.8946999999999
$ws.Range("N61").Value = -1459.6666
$ws.Range("H74").Value = 13875.74
$ws.Range("I74").Value = 17733.729
$ws.Range("J74").Value = 1230.1111
$ws.Range("K74").Value = 17733.729
$ws.Range("L74").Value = 1230.1111
$ws.Range("M74").Value = -16859.729
$ws.Range("N74").Value = -2978.1111
$ws.Range("H77").Value = 13875.74
$ws.Range("I77").Value = 17733.729
$ws.Range("J77").Value = 1230.1111
$ws.Range("K77").Value = 88668.64499999999
$ws.Range("L77").Value = 6150.5555
$ws.Range("M77").Value = -84300.64499999999
$ws.Range("N77").Value = -14886.5555
$ws.Range("H109").Value = 34984.668
$ws.Range("J109").Value = 34984.668
$ws.Range("L109").Value = 34984.668
$ws.Range("N109").Value = -37758.668
$ws.Range("H136").Value = 872.63635
$ws.Range("I136").Value = 846.8946999999999
$ws.Range("J136").Value = 1035.6666
$ws.Range("K136").Value = 2540.6841
$ws.Range("L136").Value = 3106.9998
$ws.Range("M136").Value = 9.315900000000056
$ws.Range("N136").Value = -8206.9998

# --- Sheet: BSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 29500
$ws.Range("J108").Value = 29500
$ws.Range("L108").Value = 29500
$ws.Range("N108").Value = -37180
$ws.Range("H134").Value = 22554.885
$ws.Range("I134").Value = 1182.1428
$ws.Range("J134").Value = 112320.4
$ws.Range("K134").Value = 3546.4284
$ws.Range("L134").Value = 336961.2
$ws.Range("M134").Value = -1011.4284
$ws.Range("N134").Value = -342031.2

# --- Sheet: CRP (39 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18167.146
$ws.Range("I31").Value = 19140.559
$ws.Range("J31").Value = 13439.143
$ws.Range("K31").Value = 19140.559
$ws.Range("L31").Value = 13439.143
$ws.Range("M31").Value = -18845.559
$ws.Range("N31").Value = -14029.143
$ws.Range("H34").Value = 18167.146
$ws.Range("I34").Value = 19140.559
$ws.Range("J34").Value = 13439.143
$ws.Range("K34").Value = 19140.559
$ws.Range("L34").Value = 13439.143
$ws.Range("M34").Value = -18938.559
$ws.Range("N34").Value = -13843.143
$ws.Range("H58").Value = 979.3333
$ws.Range("I58").Value = 672.1212
$ws.Range("J58").Value = 1401.75
$ws.Range("K58").Value = 672.1212
$ws.Range("L58").Value = 1401.75
$ws.Range("M58").Value = -469.1212
$ws.Range("N58").Value = -1807.75
$ws.Range("H60").Value = 9499.25
$ws.Range("J60").Value = 9499.25
$ws.Range("L60").Value = 9499.25
$ws.Range("N60").Value = -10521.25
$ws.Range("H134").Value = 985.2679000000001
$ws.Range("I134").Value = 952.9143
$ws.Range("J134").Value = 1039.1904
$ws.Range("K134").Value = 2858.7429
$ws.Range("L134").Value = 3117.5712
$ws.Range("M134").Value = -323.7429000000002
$ws.Range("N134").Value = -8187.5712
$ws.Range("H136").Value = 979.3333
$ws.Range("I136").Value = 672.1212
$ws.Range("J136").Value = 1401.75
$ws.Range("K136").Value = 2016.3636
$ws.Range("L136").Value = 4205.25
$ws.Range("M136").Value = 533.6363999999999
$ws.Range("N136").Value = -9305.25

# --- Sheet: CUL (36 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 27.666666
$ws.Range("I2").Value = 47
$ws.Range("J2").Value = 18
$ws.Range("K2").Value = 282
$ws.Range("L2").Value = 108
$ws.Range("M2").Value = -169
$ws.Range("N2").Value = -334
$ws.Range("H11").Value = 72143096
$ws.Range("I11").Value = 84166910
$ws.Range("J11").Value = 199.5
$ws.Range("K11").Value = 252500730
$ws.Range("L11").Value = 598.5
$ws.Range("M11").Value = -252500590
$ws.Range("N11").Value = -878.5
$ws.Range("H26").Value = 398.5
$ws.Range("I26").Value = 195
$ws.Range("K26").Value = 585
$ws.Range("M26").Value = -297
$ws.Range("H34").Value = 13889463
$ws.Range("J34").Value = 13889463
$ws.Range("L34").Value = 41668389
$ws.Range("N34").Value = -41668557
$ws.Range("H131").Value = 69445070
$ws.Range("I131").Value = 436.125
$ws.Range("J131").Value = 125000780
$ws.Range("K131").Value = 1308.375
$ws.Range("L131").Value = 375002340
$ws.Range("M131").Value = 3731.625
$ws.Range("N131").Value = -375012420
$ws.Range("H136").Value = 13891214
$ws.Range("I136").Value = 2183.75
$ws.Range("J136").Value = 20835730
$ws.Range("K136").Value = 6551.25
$ws.Range("L136").Value = 62507190
$ws.Range("M136").Value = -1451.25
$ws.Range("N136").Value = -62517390

# --- Sheet: GSM (14 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1066
$ws.Range("I122").Value = 980.6
$ws.Range("J122").Value = 1188
$ws.Range("K122").Value = 2941.8
$ws.Range("L122").Value = 3564
$ws.Range("M122").Value = -491.8000000000002
$ws.Range("N122").Value = -8464
$ws.Range("H132").Value = 17069.793
$ws.Range("I132").Value = 1066.6041
$ws.Range("J132").Value = 68280
$ws.Range("K132").Value = 3199.8123
$ws.Range("L132").Value = 204840
$ws.Range("M132").Value = -669.8123000000001
$ws.Range("N132").Value = -209900

# --- Sheet: LTW (28 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2629.9
$ws.Range("I7").Value = 1682.8889
$ws.Range("J7").Value = 4596.769
$ws.Range("K7").Value = 1682.8889
$ws.Range("L7").Value = 4596.769
$ws.Range("M7").Value = -1570.8889
$ws.Range("N7").Value = -4820.769
$ws.Range("H122").Value = 3085.875
$ws.Range("I122").Value = 3125.1
$ws.Range("J122").Value = 2889.75
$ws.Range("K122").Value = 9375.299999999999
$ws.Range("L122").Value = 8669.25
$ws.Range("M122").Value = -6925.299999999999
$ws.Range("N122").Value = -13569.25
$ws.Range("H126").Value = 2629.9
$ws.Range("I126").Value = 1682.8889
$ws.Range("J126").Value = 4596.769
$ws.Range("K126").Value = 5048.6667
$ws.Range("L126").Value = 13790.307
$ws.Range("M126").Value = -2578.6667
$ws.Range("N126").Value = -18730.307
$ws.Range("H132").Value = 166646.1
$ws.Range("I132").Value = 34702.35
$ws.Range("J132").Value = 775617.25
$ws.Range("K132").Value = 104107.05
$ws.Range("L132").Value = 2326851.75
$ws.Range("M132").Value = -101577.05
$ws.Range("N132").Value = -2331911.75

# --- Sheet: WVR (25 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1113.7
$ws.Range("I81").Value = 730.6667
$ws.Range("J81").Value = 1688.25
$ws.Range("K81").Value = 1461.3334
$ws.Range("L81").Value = 3376.5
$ws.Range("M81").Value = -400.3334
$ws.Range("N81").Value = -5498.5
$ws.Range("H84").Value = 1113.7
$ws.Range("I84").Value = 730.6667
$ws.Range("J84").Value = 1688.25
$ws.Range("K84").Value = 7306.666999999999
$ws.Range("L84").Value = 16882.5
$ws.Range("M84").Value = -2002.666999999999
$ws.Range("N84").Value = -27490.5
$ws.Range("H132").Value = 2311.7742
$ws.Range("I132").Value = 600.2093
$ws.Range("K132").Value = 1800.6279
$ws.Range("M132").Value = 729.3721
$ws.Range("H136").Value = 1006189.6
$ws.Range("I136").Value = 1171865.2
$ws.Range("J136").Value = 500878.94
$ws.Range("K136").Value = 3515595.6
$ws.Range("L136").Value = 1502636.82
$ws.Range("M136").Value = -3513045.6
$ws.Range("N136").Value = -1507736.82
